$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "24.972.66"
$ws.Range("E2").Value = "  -3.93%  "
$ws.Range("D3").Value = "1.640.82"
$ws.Range("E3").Value = "  -5.82%  "
$ws.Range("E4").Value = "  -0.14%  "
$ws.Range("D5").Value = "232.94"
$ws.Range("E5").Value = "  -5.79%  "
$ws.Range("D6").Value = "0.9999"
$ws.Range("E6").Value = "  -0.05%  "
$ws.Range("D7").Value = "0.4762"
$ws.Range("E7").Value = "  -5.23%  "
$ws.Range("B8").Value = "Cardano"
$ws.Range("C8").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D8").Value = "0.2585"
$ws.Range("E8").Value = "  -5.71%  "
$ws.Range("B9").Value = "Dogecoin"
$ws.Range("C9").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D9").Value = "0.06082"
$ws.Range("E9").Value = "  -1.88%  "
$ws.Range("B10").Value = "TRON"
$ws.Range("C10").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D10").Value = "0.07035"
$ws.Range("E10").Value = "  -3.23%  "
$ws.Range("B11").Value = "WrappedEther"
$ws.Range("C11").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D11").Value = "1.647.09"
$ws.Range("E11").Value = "  -5.44%  "
$ws.Range("B12").Value = "Solana"
$ws.Range("C12").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D12").Value = "14.56"
$ws.Range("E12").Value = "  -4.26%  "
$ws.Range("B13").Value = "Polygon"
$ws.Range("C13").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D13").Value = "0.5885"
$ws.Range("E13").Value = "  -10.14%  "
$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").Value = "4.326"
$ws.Range("E14").Value = "  -8.53%  "
$ws.Range("B15").Value = "Litecoin"
$ws.Range("C15").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D15").Value = "73.66"
$ws.Range("E15").Value = "  -5.45%  "
$ws.Range("B16").Value = "Dai"
$ws.Range("C16").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D16").Value = "0.9999"
$ws.Range("E16").Value = "  -0.05%  "
$ws.Range("B17").Value = "BinanceUSD"
$ws.Range("C17").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D17").Value = "0.9998"
$ws.Range("E17").Value = "  +0.00%  "
$ws.Range("B18").Value = "WrappedBTC"
$ws.Range("C18").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D18").Value = "24.959.67"
$ws.Range("E18").Value = "  -4.02%  "
$ws.Range("B19").Value = "ShibaInu"
$ws.Range("C19").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D19").Value = "0.000006578"
$ws.Range("E19").Value = "  -4.16%  "
$ws.Range("B20").Value = "Avalanche"
$ws.Range("C20").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D20").Value = "11.23"
$ws.Range("E20").Value = "  -5.67%  "
$ws.Range("B21").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C21").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D21").Value = "1.858.22"
$ws.Range("E21").Value = "  -5.52%  "
$ws.Range("B22").Value = "Uniswap"
$ws.Range("C22").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D22").Value = "4.308"
$ws.Range("E22").Value = "  -6.83%  "
$ws.Range("B23").Value = "Cosmos"
$ws.Range("C23").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D23").Value = "8.543"
$ws.Range("E23").Value = "  -2.62%  "
$ws.Range("B24").Value = "Chainlink"
$ws.Range("C24").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D24").Value = "5.226"
$ws.Range("E24").Value = "  -3.44%  "
$ws.Range("B25").Value = "Monero"
$ws.Range("C25").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D25").Value = "133.78"
$ws.Range("E25").Value = "  -0.61%  "
$ws.Range("B26").Value = "EthereumClassic"
$ws.Range("C26").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D26").Value = "14.89"
$ws.Range("E26").Value = "  -2.87%  "
$ws.Range("B27").Value = "Toncoin"
$ws.Range("C27").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D27").Value = "1.380"
$ws.Range("E27").Value = "  -7.95%  "
$ws.Range("B28").Value = "BitcoinCash"
$ws.Range("C28").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D28").Value = "104.15"
$ws.Range("E28").Value = "  -1.68%  "
$ws.Range("B29").Value = "LidoDAOToken"
$ws.Range("C29").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D29").Value = "1.639"
$ws.Range("E29").Value = "  -8.70%  "
$ws.Range("B30").Value = "InternetComputer(DFINITY)"
$ws.Range("C30").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D30").Value = "3.887"
$ws.Range("E30").Value = "  -2.82%  "
$ws.Range("B31").Value = "Filecoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D31").Value = "3.572"
$ws.Range("E31").Value = "  -3.66%  "
$ws.Range("D32").Value = "0.07582"
$ws.Range("E32").Value = "  -7.19%  "
$ws.Range("B33").Value = "Frax"
$ws.Range("C33").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D33").Value = "0.9991"
$ws.Range("E33").Value = "  -0.04%  "
$ws.Range("B34").Value = "Hedera"
$ws.Range("C34").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D34").Value = "0.04271"
$ws.Range("E34").Value = "  -10.16%  "
$ws.Range("B35").Value = "HuobiToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D35").Value = "2.570"
$ws.Range("E35").Value = "  -3.67%  "
$ws.Range("B36").Value = "ImmutableX"
$ws.Range("C36").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D36").Value = "0.5938"
$ws.Range("E36").Value = "  -3.51%  "
$ws.Range("B37").Value = "ARBITRUM"
$ws.Range("C37").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D37").Value = "0.9292"
$ws.Range("E37").Value = "  -7.04%  "
$ws.Range("B38").Value = "MXToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D38").Value = "2.574"
$ws.Range("E38").Value = "  -6.59%  "
$ws.Range("B39").Value = "TrustWalletToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D39").Value = "0.8696"
$ws.Range("E39").Value = "  +8.52%  "
$ws.Range("B40").Value = "PaxDollar"
$ws.Range("C40").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D40").Value = "0.9995"
$ws.Range("E40").Value = "  -0.07%  "
$ws.Range("B41").Value = "PaxosStandard"
$ws.Range("C41").Value = "https://coinranking.com/coin/B8xT718SbVhhh+paxosstandard-pax"
$ws.Range("D41").Value = "0.9996"
$ws.Range("E41").Value = "  -0.15%  "
$ws.Range("D42").Value = "0.01496"
$ws.Range("E42").Value = "  -7.88%  "
$ws.Range("D43").Value = "98.63"
$ws.Range("E43").Value = "  -2.59%  "
$ws.Range("D44").Value = "1.768"
$ws.Range("E44").Value = "  -8.86%  "
$ws.Range("D45").Value = "0.3709"
$ws.Range("E45").Value = "  -5.58%  "
$ws.Range("E46").Value = "  -7.00%  "
$ws.Range("D47").Value = "0.1100"
$ws.Range("E47").Value = "  -6.66%  "
$ws.Range("D48").Value = "6.104"
$ws.Range("E48").Value = "  -4.77%  "
$ws.Range("D49").Value = "0.05200"
$ws.Range("E49").Value = "  -1.79%  "
$ws.Range("E50").Value = "  -0.12%  "
$ws.Range("D51").Value = "28.68"
$ws.Range("E51").Value = "  -7.46%  "
